# "perform testcase with excel"
#
# The workbook originally has two sheets ("bus", "cab"). This adds a
# third sheet, "Test_suite", containing a small TC_ID / RUNMODE test
# matrix, and leaves it as the active sheet/tab when the file is saved.

$wb = $excel.ActiveWorkbook

# Re-assert the existing selection on "bus" (A8) before we touch anything
# else, so adding/activating the new sheet later doesn't leave "bus"
# marked as the selected tab.
$busSheet = $wb.Worksheets.Item("bus")
[void]$busSheet.Range("A8").Select()

# Add the new worksheet after the last existing sheet ("cab") and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$testSheet = $wb.Worksheets.Add($null, $lastSheet)
$testSheet.Name = "Test_suite"

# Header row.
$testSheet.Range("A1").Value = "TC_ID"
$testSheet.Range("B1").Value = "RUNMODE"

# Test matrix rows.
$testSheet.Range("A2").Value = "Bus"
$testSheet.Range("B2").Value = "Y"

$testSheet.Range("A3").Value = "Cab"
$testSheet.Range("B3").Value = "Y"

$testSheet.Range("A4").Value = "Hotel"
$testSheet.Range("B4").Value = "N"

# Leave the new sheet active with B4 selected, matching the saved view.
[void]$testSheet.Range("B4").Select()
